$d = $word.ActiveDocument

$replacements = @(
    @{old="143÷8="; new="950÷2="},
    @{old="808÷3="; new="816÷5="},
    @{old="342÷8="; new="145÷2="},
    @{old="308÷5="; new="937÷2="},
    @{old="140÷6="; new="189÷2="},
    @{old="990÷9="; new="243÷6="},
    @{old="784÷5="; new="181÷7="},
    @{old="662÷9="; new="860÷8="},
    @{old="947÷8="; new="202÷4="},
    @{old="214÷8="; new="102÷6="},
    @{old="236÷8="; new="466÷4="},
    @{old="508÷2="; new="637÷8="},
    @{old="899÷4="; new="246÷8="},
    @{old="482÷6="; new="779÷8="},
    @{old="437÷7="; new="910÷9="},
    @{old="291÷3="; new="162÷4="},
    @{old="299÷4="; new="900÷6="},
    @{old="667÷4="; new="742÷4="},
    @{old="292÷6="; new="596÷9="},
    @{old="203÷2="; new="900÷7="},
    @{old="882÷8="; new="132÷4="},
    @{old="731÷7="; new="573÷4="},
    @{old="589÷7="; new="844÷4="},
    @{old="201÷7="; new="668÷6="},
    @{old="957÷6="; new="664÷5="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
